$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 905.6667
$ws.Range("J17").Value = 849.6923
$ws.Range("L17").Value = 2549.0769
$ws.Range("N17").Value = -2885.0769

$ws.Range("H19").Value = 562.875
$ws.Range("J19").Value = 571.8570999999999
$ws.Range("L19").Value = 571.8570999999999
$ws.Range("N19").Value = -921.8570999999999

$ws.Range("H126").Value = 59866.2
$ws.Range("J126").Value = 59866.2
$ws.Range("L126").Value = 59866.2
$ws.Range("N126").Value = -69746.2

$ws.Range("H130").Value = 33997.8
$ws.Range("J130").Value = 33997.8
$ws.Range("L130").Value = 33997.8
$ws.Range("N130").Value = -44037.8

$ws.Range("H137").Value = 2600.4
$ws.Range("I137").Value = 2334.125
$ws.Range("J137").Value = 3665.5
$ws.Range("K137").Value = 7002.375
$ws.Range("L137").Value = 10996.5
$ws.Range("M137").Value = -4452.375
$ws.Range("N137").Value = -16096.5

$ws.Range("H138").Value = 2577.0444
$ws.Range("I138").Value = 1648.2174
$ws.Range("J138").Value = 3548.0908
$ws.Range("K138").Value = 4944.6522
$ws.Range("L138").Value = 10644.2724
$ws.Range("M138").Value = 195.3477999999996
$ws.Range("N138").Value = -20924.2724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 210.66667
$ws.Range("I4").Value = 210.66667
$ws.Range("K4").Value = 210.66667
$ws.Range("M4").Value = -94.66667000000001

$ws.Range("H5").Value = 555000000
$ws.Range("I5").Value = 110000000
$ws.Range("J5").Value = 1000000000
$ws.Range("K5").Value = 110000000
$ws.Range("L5").Value = 1000000000
$ws.Range("M5").Value = -109999888
$ws.Range("N5").Value = -1000000224

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H92").Value = 50399.1
$ws.Range("J92").Value = 50399.1
$ws.Range("L92").Value = 50399.1
$ws.Range("N92").Value = -55391.1

$ws.Range("H132").Value = 3753.2104
$ws.Range("I132").Value = 3835.6875
$ws.Range("K132").Value = 11507.0625
$ws.Range("M132").Value = -8977.0625

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 555000000
$ws.Range("I4").Value = 110000000
$ws.Range("J4").Value = 1000000000
$ws.Range("K4").Value = 110000000
$ws.Range("L4").Value = 1000000000
$ws.Range("M4").Value = -109999885
$ws.Range("N4").Value = -1000000230

$ws.Range("H92").Value = 130394
$ws.Range("J92").Value = 130394
$ws.Range("L92").Value = 130394
$ws.Range("N92").Value = -135386

$ws.Range("H99").Value = 103752.5
$ws.Range("I99").Value = 400010
$ws.Range("K99").Value = 400010
$ws.Range("M99").Value = -398512

$ws.Range("H100").Value = 16416.334
$ws.Range("J100").Value = 16416.334
$ws.Range("L100").Value = 16416.334
$ws.Range("N100").Value = -18580.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1000000000
$ws.Range("I4").Value = 1000000000
$ws.Range("K4").Value = 1000000000
$ws.Range("M4").Value = -999999888

$ws.Range("H31").Value = 1930
$ws.Range("I31").Value = 1811.9615
$ws.Range("K31").Value = 1811.9615
$ws.Range("M31").Value = -1516.9615

$ws.Range("H34").Value = 1930
$ws.Range("I34").Value = 1811.9615
$ws.Range("K34").Value = 1811.9615
$ws.Range("M34").Value = -1609.9615

$ws.Range("H59").Value = 96516.14
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 96516.14
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 96516.14
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -98806.14

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 599.36365
$ws.Range("I29").Value = 769.4286
$ws.Range("J29").Value = 301.75
$ws.Range("K29").Value = 2308.2858
$ws.Range("L29").Value = 905.25
$ws.Range("M29").Value = -2031.2858
$ws.Range("N29").Value = -1459.25

$ws.Range("H34").Value = 2586
$ws.Range("I34").Value = 179
$ws.Range("J34").Value = 2929.8572
$ws.Range("K34").Value = 537
$ws.Range("L34").Value = 8789.571599999999
$ws.Range("M34").Value = -453
$ws.Range("N34").Value = -8957.571599999999

$ws.Range("H44").Value = 94
$ws.Range("I44").Value = 94
$ws.Range("K44").Value = 282
$ws.Range("M44").Value = 116

$ws.Range("H46").Value = 100452.5
$ws.Range("I46").Value = 333460
$ws.Range("K46").Value = 1000380
$ws.Range("M46").Value = -1000289

$ws.Range("H68").Value = 1959.9231
$ws.Range("I68").Value = 999.2
$ws.Range("J68").Value = 2560.375
$ws.Range("K68").Value = 2997.6
$ws.Range("L68").Value = 7681.125
$ws.Range("M68").Value = -2186.6
$ws.Range("N68").Value = -9303.125

$ws.Range("H71").Value = 1959.9231
$ws.Range("I71").Value = 999.2
$ws.Range("J71").Value = 2560.375
$ws.Range("K71").Value = 8992.800000000001
$ws.Range("L71").Value = 23043.375
$ws.Range("M71").Value = -4936.800000000001
$ws.Range("N71").Value = -31155.375

$ws.Range("H95").Value = 4930
$ws.Range("I95").Value = 4890
$ws.Range("J95").Value = 4950
$ws.Range("K95").Value = 14670
$ws.Range("L95").Value = 14850
$ws.Range("M95").Value = -12611
$ws.Range("N95").Value = -18968

$ws.Range("H121").Value = 1316.1666
$ws.Range("I121").Value = 599.25
$ws.Range("J121").Value = 2750
$ws.Range("K121").Value = 1797.75
$ws.Range("L121").Value = 8250
$ws.Range("M121").Value = -487.75
$ws.Range("N121").Value = -10870

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3519.3547
$ws.Range("I102").Value = 2712.7083
$ws.Range("J102").Value = 6285
$ws.Range("K102").Value = 2712.7083
$ws.Range("L102").Value = 6285
$ws.Range("M102").Value = -1090.7083
$ws.Range("N102").Value = -9529

$ws.Range("H126").Value = 3501.5
$ws.Range("I126").Value = 3005.5
$ws.Range("K126").Value = 9016.5
$ws.Range("M126").Value = -6546.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2394.75
$ws.Range("I7").Value = 2348.4
$ws.Range("J7").Value = 2427.8572
$ws.Range("K7").Value = 2348.4
$ws.Range("L7").Value = 2427.8572
$ws.Range("M7").Value = -2236.4
$ws.Range("N7").Value = -2651.8572

$ws.Range("H93").Value = 2498.0667
$ws.Range("I93").Value = 2413
$ws.Range("K93").Value = 2413
$ws.Range("M93").Value = -1165

$ws.Range("H106").Value = 26465
$ws.Range("J106").Value = 26465
$ws.Range("L106").Value = 26465
$ws.Range("N106").Value = -28989

$ws.Range("H108").Value = 42250
$ws.Range("J108").Value = 42250
$ws.Range("L108").Value = 42250
$ws.Range("N108").Value = -49930

$ws.Range("H122").Value = 5016.5454
$ws.Range("I122").Value = 5016.5454
$ws.Range("K122").Value = 15049.6362
$ws.Range("M122").Value = -12599.6362

$ws.Range("H126").Value = 2394.75
$ws.Range("I126").Value = 2348.4
$ws.Range("J126").Value = 2427.8572
$ws.Range("K126").Value = 7045.200000000001
$ws.Range("L126").Value = 7283.571599999999
$ws.Range("M126").Value = -4575.200000000001
$ws.Range("N126").Value = -12223.5716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2851.8
$ws.Range("I122").Value = 1919.6666
$ws.Range("K122").Value = 5758.9998
$ws.Range("M122").Value = -3308.9998

$ws.Range("H126").Value = 5010.905
$ws.Range("I126").Value = 5450.231
$ws.Range("J126").Value = 4297
$ws.Range("K126").Value = 16350.693
$ws.Range("L126").Value = 12891
$ws.Range("M126").Value = -13880.693
$ws.Range("N126").Value = -17831
